$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-20 Wednesday" "2024-03-21 Thursday"

Replace-Text "752÷8=" "411÷7="
Replace-Text "698÷4=" "341÷6="
Replace-Text "298÷3=" "178÷8="
Replace-Text "719÷6=" "697÷8="
Replace-Text "928÷8=" "642÷4="
Replace-Text "509÷4=" "878÷8="
Replace-Text "785÷6=" "766÷2="
Replace-Text "356÷3=" "841÷7="
Replace-Text "164÷7=" "795÷5="
Replace-Text "932÷7=" "916÷5="
Replace-Text "282÷9=" "964÷2="
Replace-Text "120÷6=" "932÷6="
Replace-Text "887÷4=" "672÷5="
Replace-Text "125÷3=" "496÷6="
Replace-Text "557÷9=" "738÷7="
Replace-Text "419÷4=" "210÷3="
Replace-Text "750÷5=" "751÷3="
Replace-Text "203÷4=" "305÷7="
Replace-Text "184÷3=" "495÷6="
Replace-Text "176÷7=" "338÷3="
Replace-Text "871÷5=" "948÷2="
Replace-Text "667÷5=" "750÷7="
Replace-Text "536÷5=" "322÷2="
Replace-Text "654÷8=" "923÷6="
Replace-Text "744÷6=" "176÷9="
